$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "VALOR MORA" total and "Cant. Periodos" count -------------
$ws.Range("E11").Value = 227760
$ws.Range("F13").Value = 2

# --- Give the soon-to-be-last data row (row 19) the closing bottom-border
#     formatting that row 21 (the current last row) carries, before the
#     now-obsolete rows 20:21 are removed. -------------------------------
$ws.Range("B21:J21").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)

# --- Replace worker period 2506 (row 17, Andres) with Adriana's 2507 row -
$ws.Range("C17").Value = "1101816566"
$ws.Range("D17").Value = "ADRIANA MARCELA RIVAS PEREZ"
$ws.Range("E17").Value = "2507"

# --- Replace worker period 2505 (row 18, Andres) with Andres' new 2508 --
$ws.Range("E18").Value = "2508"

# --- Replace worker period 2505 (row 19, Adriana, now reformatted) with
#     Adriana's new 2508 row ------------------------------------------
$ws.Range("E19").Value = "2508"

# --- Remove the two now-duplicate trailing rows (old 2506/2505 Adriana) -
$ws.Rows("20:21").Delete()
